$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their text formatting (avoid numeric auto-conversion)
$priceCells = @("D2","D3","D5","D6","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D38","D39","D40","D41","D42","D45","D46","D48","D50","D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated crypto values scraped on Tue Jan  9 04:56:46 UTC 2024
$ws.Range("D2").Value = "46.815.48"
$ws.Range("E2").Value = "  +6.53%  "
$ws.Range("D3").Value = "2.314.22"
$ws.Range("E3").Value = "  +5.30%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "301.04"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "103.51"
$ws.Range("E6").Value = "  +15.94%  "
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +10.15%  "
$ws.Range("D10").Value = "37.45"
$ws.Range("E10").Value = "  +16.38%  "
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  +4.59%  "
$ws.Range("D12").Value = "7.46"
$ws.Range("E12").Value = "  +9.85%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "2.665.27"
$ws.Range("E14").Value = "  +5.24%  "
$ws.Range("D15").Value = "2.304.23"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "14.08"
$ws.Range("E16").Value = "  +7.38%  "
$ws.Range("D17").Value = "0.828"
$ws.Range("E17").Value = "  +6.99%  "
$ws.Range("D18").Value = "46.796.87"
$ws.Range("E18").Value = "  +7.26%  "
$ws.Range("D19").Value = "13.37"
$ws.Range("E19").Value = "  +23.37%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +7.17%  "
$ws.Range("D21").Value = "6.19"
$ws.Range("E21").Value = "  +5.94%  "
$ws.Range("D22").Value = "67.29"
$ws.Range("E22").Value = "  +6.70%  "
$ws.Range("D23").Value = "248.99"
$ws.Range("E23").Value = "  +7.56%  "
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  +7.15%  "
$ws.Range("E25").Value = "  +8.65%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").Value = "43.59"
$ws.Range("E27").Value = "  +20.86%  "
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").Value = "9.99"
$ws.Range("E29").Value = "  +8.24%  "
$ws.Range("D30").Value = "20.27"
$ws.Range("E30").Value = "  +5.72%  "
$ws.Range("D31").Value = "5.81"
$ws.Range("E31").Value = "  +9.51%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0806"
$ws.Range("E32").Value = "  +9.39%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "147.32"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "2.63"
$ws.Range("E34").Value = "  +4.95%  "
$ws.Range("D35").Value = "3.13"
$ws.Range("E35").Value = "  +9.18%  "
$ws.Range("E36").Value = "  +9.09%  "
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  +9.43%  "
$ws.Range("D39").Value = "15.97"
$ws.Range("E39").Value = "  +21.70%  "
$ws.Range("D40").Value = "4.07"
$ws.Range("E40").Value = "  +15.01%  "
$ws.Range("D41").Value = "3.48"
$ws.Range("E41").Value = "  +12.25%  "
$ws.Range("D42").Value = "0.0309"
$ws.Range("E42").Value = "  +9.64%  "
$ws.Range("E43").Value = "  +21.64%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "1.847.81"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("D46").Value = "89.29"
$ws.Range("E46").Value = "  +23.11%  "
$ws.Range("E47").Value = "  +14.20%  "
$ws.Range("D48").Value = "75.46"
$ws.Range("E48").Value = "  +15.06%  "
$ws.Range("E49").Value = "  +9.80%  "
$ws.Range("D50").Value = "98.04"
$ws.Range("E50").Value = "  +6.62%  "
$ws.Range("D51").Value = "55.01"
$ws.Range("E51").Value = "  +12.03%  "
